# Applies the diff to the tail of "Leetcode Log.docx":
#   - strips <w:rFonts w:hint="eastAsia"/> from the paragraph-mark rPr of the
#     empty paragraph before "054_SpiralMatrix 1Y" and from the
#     "054_SpiralMatrix 1Y" paragraph itself (run-level rFonts is untouched)
#   - removes the _GoBack bookmark from the end of the
#     "Mind all possible conditions." paragraph
#   - appends two new paragraphs after the trailing empty paragraph:
#       "058_LengthOfLastWord 3Y"
#       "Forgot to consider the space could be at the end and there could be
#        multiple spaces." (ending with the relocated _GoBack bookmark)

$d = $word.ActiveDocument

$wNs = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'

# Locate the anchor paragraphs by their text, walking from the end of the
# document so the edit is robust to anything earlier in the log.
$total = $d.Paragraphs.Count

$idxSpiral = -1
$idxConditions = -1
for ($i = $total; $i -ge 1; $i--) {
    $t = $d.Paragraphs.Item($i).Range.Text.TrimEnd([char]13, [char]7)
    if ($idxConditions -eq -1 -and $t -eq "Mind all possible conditions.") {
        $idxConditions = $i
    }
    if ($idxSpiral -eq -1 -and $t -eq "054_SpiralMatrix 1Y") {
        $idxSpiral = $i
    }
    if ($idxSpiral -ne -1 -and $idxConditions -ne -1) {
        break
    }
}

if ($idxSpiral -eq -1) {
    throw "Could not locate the '054_SpiralMatrix 1Y' paragraph."
}
if ($idxConditions -eq -1) {
    throw "Could not locate the 'Mind all possible conditions.' paragraph."
}

$idxEmptyBeforeSpiral = $idxSpiral - 1
$idxLast = $d.Paragraphs.Count

# 1) Empty paragraph right before "054_SpiralMatrix 1Y": drop the
#    rFonts hint="eastAsia" from the paragraph mark's rPr.
$pEmpty = $d.Paragraphs.Item($idxEmptyBeforeSpiral)
$xmlEmpty = "<w:p $wNs><w:pPr><w:rPr><w:color w:val=""000000"" w:themeColor=""text1""/><w:lang w:eastAsia=""zh-CN""/></w:rPr></w:pPr></w:p>"
[void]$pEmpty.Range.InsertXML($xmlEmpty)

# 2) "054_SpiralMatrix 1Y" paragraph: drop the rFonts hint from the
#    paragraph mark's rPr, but keep it on the run itself.
$pSpiral = $d.Paragraphs.Item($idxSpiral)
$xmlSpiral = "<w:p $wNs><w:pPr><w:rPr><w:color w:val=""000000"" w:themeColor=""text1""/><w:lang w:eastAsia=""zh-CN""/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:hint=""eastAsia""/><w:color w:val=""000000"" w:themeColor=""text1""/><w:lang w:eastAsia=""zh-CN""/></w:rPr><w:t>054_SpiralMatrix 1Y</w:t></w:r></w:p>"
[void]$pSpiral.Range.InsertXML($xmlSpiral)

# 3) "Mind all possible conditions." paragraph: remove the trailing
#    _GoBack bookmark (bookmarkStart/bookmarkEnd) while keeping the two runs.
$pConditions = $d.Paragraphs.Item($idxConditions)
$xmlConditions = "<w:p $wNs><w:pPr><w:rPr><w:color w:val=""000000"" w:themeColor=""text1""/><w:lang w:eastAsia=""zh-CN""/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:hint=""eastAsia""/><w:color w:val=""000000"" w:themeColor=""text1""/><w:lang w:eastAsia=""zh-CN""/></w:rPr><w:t>Mind all possible</w:t></w:r><w:r><w:rPr><w:color w:val=""000000"" w:themeColor=""text1""/><w:lang w:eastAsia=""zh-CN""/></w:rPr><w:t xml:space=""preserve""> conditions.</w:t></w:r></w:p>"
[void]$pConditions.Range.InsertXML($xmlConditions)

# 4) Append two new paragraphs after the trailing empty paragraph, pushing
#    the _GoBack bookmark onto the new final paragraph.
$pLast = $d.Paragraphs.Item($idxLast)
[void]$pLast.Range.InsertParagraphAfter()

$pNew1 = $d.Paragraphs.Item($idxLast + 1)
$xmlNew1 = "<w:p $wNs><w:pPr><w:rPr><w:rFonts w:hint=""eastAsia""/><w:color w:val=""000000"" w:themeColor=""text1""/><w:lang w:eastAsia=""zh-CN""/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:hint=""eastAsia""/><w:color w:val=""000000"" w:themeColor=""text1""/><w:lang w:eastAsia=""zh-CN""/></w:rPr><w:t>058_LengthOfLastWord 3Y</w:t></w:r></w:p>"
[void]$pNew1.Range.InsertXML($xmlNew1)

$pNew1Again = $d.Paragraphs.Item($idxLast + 1)
[void]$pNew1Again.Range.InsertParagraphAfter()

$pNew2 = $d.Paragraphs.Item($idxLast + 2)
$xmlNew2 = "<w:p $wNs><w:pPr><w:rPr><w:rFonts w:hint=""eastAsia""/><w:color w:val=""000000"" w:themeColor=""text1""/><w:lang w:eastAsia=""zh-CN""/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:hint=""eastAsia""/><w:color w:val=""000000"" w:themeColor=""text1""/><w:lang w:eastAsia=""zh-CN""/></w:rPr><w:t>Forgot to consider the space could be at the end and there could be multiple spaces.</w:t></w:r><w:bookmarkStart w:id=""0"" w:name=""_GoBack""/><w:bookmarkEnd w:id=""0""/></w:p>"
[void]$pNew2.Range.InsertXML($xmlNew2)

Write-Output "Done. Paragraphs now: $($d.Paragraphs.Count)"
